$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source data added a new weekly price record for "Coco" (Vega Central
# Mapocho de Santiago) dated 2021-11-29, which belongs right after the
# header/above rows in the existing date-ordered block. It lands at row 23,
# pushing the previous rows 23-36 down to 24-37.
$ws.Rows.Item(23).Insert()

# Populate the newly inserted row 23 with the new record. The static
# descriptive columns (A,B,C,E,F,G,H,I,J,K,L,Q,R,T) are identical to every
# other row in this "Coco" block.
$ws.Cells.Item(23, 1).Value = 9
$ws.Cells.Item(23, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(23, 3).Value = "Metropolitana"
$ws.Cells.Item(23, 4).Value = 44529
$ws.Cells.Item(23, 5).Value = 13
$ws.Cells.Item(23, 6).Value = "Fruta"
$ws.Cells.Item(23, 7).Value = 100108
$ws.Cells.Item(23, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(23, 9).Value = 100108007
$ws.Cells.Item(23, 10).Value = "Coco"
$ws.Cells.Item(23, 11).Value = "Sin especificar"
$ws.Cells.Item(23, 12).Value = "Primera"
$ws.Cells.Item(23, 13).Value = 34
$ws.Cells.Item(23, 14).Value = 28000
$ws.Cells.Item(23, 15).Value = 28000
$ws.Cells.Item(23, 16).Value = 28000
$ws.Cells.Item(23, 17).Value = "$/malla 20 unidades"
$ws.Cells.Item(23, 18).Value = "Perú"
$ws.Cells.Item(23, 19).Value = 1400
$ws.Cells.Item(23, 20).Value = 20
